# Github Commands.docx - add a new run containing " .." right after the
# run/paragraph that reads "Add all new files" (list item bullet).
#
# The new run must keep exactly the same visible formatting as the run
# it follows (Segoe UI / Times New Roman(eastAsia) / Segoe UI(cs), color
# 24292E, size 7.5pt, lang eastAsia=de-DE) but be serialized as its own
# separate <w:r> element rather than being merged into the preceding run.

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Add all new files") {

        $para = $p.Range

        # Position right before the paragraph mark (end of the visible text).
        $insertStart = $para.End - 1

        # Insert the new literal text, preserving the leading space.
        $ip = $d.Range($insertStart, $insertStart)
        $ip.InsertAfter(" ..")

        # Range covering exactly the text we just inserted.
        $newRun = $d.Range($insertStart, $insertStart + 3)

        # Toggling a character property and then restoring its original
        # value forces the interop layer to keep this text as an
        # independent run (rather than silently re-merging it with the
        # identically-formatted run before it), while leaving the final
        # formatting unchanged.
        $newRun.Bold = 1
        $newRun.Bold = 0

        break
    }
}
